# Scheduled-runner price refresh: update computed columns H:N (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across all 8 job sheets with freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 622.8333
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 1146.6666
$ws.Range("K17").Value = 297
$ws.Range("L17").Value = 3439.9998
$ws.Range("M17").Value = -129
$ws.Range("N17").Value = -3775.9998
# Row 38
$ws.Range("H38").Value = 2106.2778
$ws.Range("I38").Value = 798.5
$ws.Range("J38").Value = 4721.8335
$ws.Range("K38").Value = 2395.5
$ws.Range("L38").Value = 14165.5005
$ws.Range("M38").Value = -2023.5
$ws.Range("N38").Value = -14909.5005
# Row 58
$ws.Range("H58").Value = 9286.462
$ws.Range("I58").Value = 191.55556
$ws.Range("K58").Value = 574.66668
$ws.Range("M58").Value = -424.66668
# Row 64
$ws.Range("H64").Value = 4849.5
$ws.Range("I64").Value = 4700
$ws.Range("J64").ClearContents()
$ws.Range("K64").Value = 4700
$ws.Range("L64").ClearContents()
$ws.Range("M64").Value = -4452
$ws.Range("N64").Value = -5495
# Row 67
$ws.Range("H67").Value = 4849.5
$ws.Range("I67").Value = 4700
$ws.Range("J67").ClearContents()
$ws.Range("K67").Value = 4700
$ws.Range("L67").ClearContents()
$ws.Range("M67").Value = -3842
$ws.Range("N67").Value = -6715
# Row 107
$ws.Range("H107").Value = 916.8
$ws.Range("I107").Value = 916.8
$ws.Range("K107").Value = 916.8
$ws.Range("M107").Value = 1003.2
# Row 132
$ws.Range("H132").Value = 1701.0322
$ws.Range("I132").Value = 1580.5769
$ws.Range("K132").Value = 4741.7307
$ws.Range("M132").Value = -2211.7307
# Row 137
$ws.Range("H137").Value = 55558664
$ws.Range("I137").Value = 62503372
$ws.Range("K137").Value = 187510116
$ws.Range("M137").Value = -187507566
# Row 138
$ws.Range("H138").Value = 2236.5227
$ws.Range("J138").Value = 2298.5144
$ws.Range("L138").Value = 6895.5432
$ws.Range("N138").Value = -17175.5432
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3160.7021
$ws.Range("I32").Value = 1335.0555
$ws.Range("K32").Value = 1335.0555
$ws.Range("M32").Value = -1048.0555
# Row 122
$ws.Range("H122").Value = 12298.4
$ws.Range("I122").Value = 16497.666
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 49492.99800000001
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -47042.99800000001
$ws.Range("N122").Value = -22898.5
# Row 132
$ws.Range("H132").Value = 13482.182
$ws.Range("I132").Value = 10830.7
$ws.Range("J132").Value = 39997
$ws.Range("K132").Value = 32492.1
$ws.Range("L132").Value = 119991
$ws.Range("M132").Value = -29962.1
$ws.Range("N132").Value = -125051
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2261.0417
$ws.Range("I16").Value = 2505.4211
$ws.Range("K16").Value = 2505.4211
$ws.Range("M16").Value = -2218.4211
# Row 31
$ws.Range("H31").Value = 3164.6667
$ws.Range("I31").Value = 2747.5715
$ws.Range("K31").Value = 2747.5715
$ws.Range("M31").Value = -2452.5715
# Row 34
$ws.Range("H34").Value = 3164.6667
$ws.Range("I34").Value = 2747.5715
$ws.Range("K34").Value = 2747.5715
$ws.Range("M34").Value = -2545.5715
# Row 103
$ws.Range("H103").Value = 24000
$ws.Range("I103").Value = 24000
$ws.Range("K103").Value = 24000
$ws.Range("M103").Value = -22828
# Row 105
$ws.Range("H105").Value = 6641.8823
$ws.Range("I105").Value = 8608.959999999999
$ws.Range("J105").Value = 1177.7778
$ws.Range("K105").Value = 8608.959999999999
$ws.Range("L105").Value = 1177.7778
$ws.Range("M105").Value = -6861.959999999999
$ws.Range("N105").Value = -4671.7778
# Row 113
$ws.Range("H113").Value = 2261.0417
$ws.Range("I113").Value = 2505.4211
$ws.Range("K113").Value = 2505.4211
$ws.Range("M113").Value = -335.4211
# Row 132
$ws.Range("H132").Value = 3546.6428
$ws.Range("I132").Value = 3332.24
$ws.Range("K132").Value = 9996.719999999999
$ws.Range("M132").Value = -7466.719999999999
# Row 134
$ws.Range("H134").Value = 2525.577
$ws.Range("I134").Value = 2338.8948
$ws.Range("J134").Value = 3032.2856
$ws.Range("K134").Value = 7016.6844
$ws.Range("L134").Value = 9096.856800000001
$ws.Range("M134").Value = -4481.6844
$ws.Range("N134").Value = -14166.8568
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 250000340
$ws.Range("I4").Value = 333333440
$ws.Range("J4").Value = 998
$ws.Range("K4").Value = 1000000320
$ws.Range("L4").Value = 2994
$ws.Range("M4").Value = -1000000208
$ws.Range("N4").Value = -3218
# Row 5
$ws.Range("H5").Value = 470
$ws.Range("I5").Value = 291.66666
$ws.Range("K5").Value = 874.9999799999999
$ws.Range("M5").Value = -762.9999799999999
# Row 6
$ws.Range("H6").Value = 43458.75
$ws.Range("I6").Value = 71788.57000000001
$ws.Range("J6").Value = 3797
$ws.Range("K6").Value = 215365.71
$ws.Range("L6").Value = 11391
$ws.Range("M6").Value = -215252.71
$ws.Range("N6").Value = -11617
# Row 7
$ws.Range("H7").Value = 340
$ws.Range("I7").Value = 309.16666
$ws.Range("K7").Value = 927.4999799999999
$ws.Range("M7").Value = -815.4999799999999
# Row 33
$ws.Range("H33").Value = 256.72223
$ws.Range("J33").Value = 278
$ws.Range("L33").Value = 1668
$ws.Range("N33").Value = -2234
# Row 37
$ws.Range("H37").Value = 139015.36
$ws.Range("J37").Value = 139015.36
$ws.Range("L37").Value = 417046.08
$ws.Range("N37").Value = -417270.08
# Row 40
$ws.Range("H40").Value = 202.26086
$ws.Range("I40").Value = 113.15
$ws.Range("J40").Value = 796.3333
$ws.Range("K40").Value = 452.6
$ws.Range("L40").Value = 3185.3332
$ws.Range("M40").Value = -383.6
$ws.Range("N40").Value = -3323.3332
# Row 107
$ws.Range("H107").Value = 537
$ws.Range("J107").Value = 563.73914
$ws.Range("L107").Value = 1691.21742
$ws.Range("N107").Value = -5531.21742
# Row 135
$ws.Range("H135").Value = 470
$ws.Range("I135").Value = 291.66666
$ws.Range("K135").Value = 2624.99994
$ws.Range("M135").Value = -89.9999399999997
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 25527.814
$ws.Range("I70").Value = 42853.855
$ws.Range("J70").Value = 6869
$ws.Range("K70").Value = 42853.855
$ws.Range("L70").Value = 6869
$ws.Range("M70").Value = -42583.855
$ws.Range("N70").Value = -7409
# Row 73
$ws.Range("H73").Value = 25527.814
$ws.Range("I73").Value = 42853.855
$ws.Range("J73").Value = 6869
$ws.Range("K73").Value = 42853.855
$ws.Range("L73").Value = 6869
$ws.Range("M73").Value = -41917.855
$ws.Range("N73").Value = -8741
# Row 80
$ws.Range("H80").Value = 8995.200000000001
$ws.Range("J80").Value = 8993.333000000001
$ws.Range("L80").Value = 8993.333000000001
$ws.Range("N80").Value = -10989.333
# Row 83
$ws.Range("H83").Value = 8995.200000000001
$ws.Range("J83").Value = 8993.333000000001
$ws.Range("L83").Value = 44966.665
$ws.Range("N83").Value = -54950.665
# Row 102
$ws.Range("H102").Value = 4206.7646
$ws.Range("I102").Value = 4456.091
$ws.Range("J102").Value = 3749.6667
$ws.Range("K102").Value = 4456.091
$ws.Range("L102").Value = 3749.6667
$ws.Range("M102").Value = -2834.091
$ws.Range("N102").Value = -6993.6667
# Row 122
$ws.Range("H122").Value = 3410.9333
$ws.Range("I122").Value = 3486.182
$ws.Range("K122").Value = 10458.546
$ws.Range("M122").Value = -8008.545999999998
# Row 126
$ws.Range("H126").Value = 4076.25
$ws.Range("I126").Value = 2136.6667
$ws.Range("J126").Value = 5240
$ws.Range("K126").Value = 6410.000100000001
$ws.Range("L126").Value = 15720
$ws.Range("M126").Value = -3940.000100000001
$ws.Range("N126").Value = -20660
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 58828570
$ws.Range("I7").Value = 142859890
$ws.Range("J7").Value = 6638.7
$ws.Range("K7").Value = 142859890
$ws.Range("L7").Value = 6638.7
$ws.Range("M7").Value = -142859778
$ws.Range("N7").Value = -6862.7
# Row 16
$ws.Range("H16").Value = 1737.8334
$ws.Range("J16").Value = 1986.6666
$ws.Range("L16").Value = 1986.6666
$ws.Range("N16").Value = -2326.6666
# Row 100
$ws.Range("H100").Value = 5466.3335
$ws.Range("I100").Value = 3200
$ws.Range("K100").Value = 3200
$ws.Range("M100").Value = -2659
# Row 126
$ws.Range("H126").Value = 58828570
$ws.Range("I126").Value = 142859890
$ws.Range("J126").Value = 6638.7
$ws.Range("K126").Value = 428579670
$ws.Range("L126").Value = 19916.1
$ws.Range("M126").Value = -428577200
$ws.Range("N126").Value = -24856.1
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 9799
$ws.Range("J6").Value = 9799
$ws.Range("L6").Value = 9799
$ws.Range("N6").Value = -10029
# Row 112
$ws.Range("H112").Value = 18997
$ws.Range("J112").Value = 18997
$ws.Range("L112").Value = 18997
$ws.Range("N112").Value = -21951
# Row 136
$ws.Range("H136").Value = 4487.5386
$ws.Range("I136").Value = 4451.579
$ws.Range("J136").Value = 4585.143
$ws.Range("K136").Value = 13354.737
$ws.Range("L136").Value = 13755.429
$ws.Range("M136").Value = -10804.737
$ws.Range("N136").Value = -18855.429

Write-Host "Applied scheduled Sheets update (233 cells across 8 sheets)"
